# Colocando header nos gráficos
# Adds a header label in column A (row 1) for each data table and removes
# the bold/bordered "header" style from the now-plain row labels below it.
# Also fixes a few accented Portuguese labels, removes the stray "Teto" row
# from the Emissoes sheet, and updates the "Custo Total" sheet header/values.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($ws, $addr, $text) {
    # Give the new header cell the same look (bold/border/centered) as the
    # other header cells in row 1, by copying B1's format onto it.
    $ws.Range($addr).Value = $text
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Clear-LabelStyle($ws, $addr) {
    # Remove the bold/border/centered style so the cell reverts to the
    # default (unstyled) look, like the rest of the data rows.
    $ws.Range($addr).ClearFormats()
}

# ---------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
# (MWMed)", "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# ---------------------------------------------------------------------
foreach ($idx in 1..4) {
    $ws = $wb.Worksheets.Item($idx)

    Set-HeaderCell $ws "A1" "Fonte/Tecnologia"

    Clear-LabelStyle $ws "A2"

    $ws.Range("A3").Value = "Gás Natural"
    Clear-LabelStyle $ws "A3"

    $ws.Range("A4").Value = "Carvão"
    Clear-LabelStyle $ws "A4"

    Clear-LabelStyle $ws "A5"

    $ws.Range("A6").Value = "Óleos Comb"
    Clear-LabelStyle $ws "A6"

    Clear-LabelStyle $ws "A7"

    $ws.Range("A8").Value = "Eólica"
    Clear-LabelStyle $ws "A8"

    Clear-LabelStyle $ws "A9"

    Clear-LabelStyle $ws "A10"

    $ws.Range("A11").Value = "Pot. Compl."
    Clear-LabelStyle $ws "A11"

    Clear-LabelStyle $ws "A12"
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "A1" "Período"

$ws5.Range("A2").Value = "P.Médio"
Clear-LabelStyle $ws5 "A2"

$ws5.Range("A3").Value = "P.Crítico"
Clear-LabelStyle $ws5 "A3"

# Remove the extra "Teto" row entirely (row 4)
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "A1" "Tipo Expansão"
# Leading apostrophe forces Excel to store the numeric-looking label as
# text (matching how the other sheets store "2015" as text, not a number).
$ws6.Range("B1").Value = "'2015"

$ws6.Range("A2").Value = "Expansão Centralizada"
Clear-LabelStyle $ws6 "A2"
$ws6.Range("B2").Value = 582

$ws6.Range("A3").Value = "Expansão por GD"
Clear-LabelStyle $ws6 "A3"
$ws6.Range("B3").Value = 99
